$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values (e.g. '102.40')
# are not coerced into floating point numbers, matching the inline-string source data.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '51.141.51'
$ws.Cells.Item(2, 5).Value = '  +0.55%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.963.12'
$ws.Cells.Item(3, 5).Value = '  +1.41%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  -0.07%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '379.88'
$ws.Cells.Item(5, 5).Value = '  +2.45%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '102.40'
$ws.Cells.Item(6, 5).Value = '  +3.05%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.545'
$ws.Cells.Item(7, 5).Value = '  +3.05%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.00%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.589'
$ws.Cells.Item(9, 5).Value = '  +2.75%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '36.51'
$ws.Cells.Item(10, 5).Value = '  +2.62%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.61%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '0.0858'
$ws.Cells.Item(12, 5).Value = '  +2.61%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '3.424.31'
$ws.Cells.Item(13, 5).Value = '  +1.35%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '7.77'
$ws.Cells.Item(14, 5).Value = '  +5.80%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '18.27'
$ws.Cells.Item(15, 5).Value = '  +2.91%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.958.62'
$ws.Cells.Item(16, 5).Value = '  +1.08%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '11.12'
$ws.Cells.Item(17, 5).Value = '  -0.65%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +4.30%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '51.162.22'
$ws.Cells.Item(19, 5).Value = '  +0.60%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '3.19'
$ws.Cells.Item(20, 5).Value = '  +3.01%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '12.54'
$ws.Cells.Item(21, 5).Value = '  +3.30%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.0₃0958'
$ws.Cells.Item(22, 5).Value = '  +1.67%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '70.09'
$ws.Cells.Item(23, 5).Value = '  +3.03%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '266.84'
$ws.Cells.Item(24, 5).Value = '  +1.93%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '3.20'
$ws.Cells.Item(25, 5).Value = '  +4.47%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'RenderToken'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(26, 4).Value = '7.50'
$ws.Cells.Item(26, 5).Value = '  +4.66%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'Filecoin'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(27, 4).Value = '7.77'
$ws.Cells.Item(27, 5).Value = '  -2.24%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -0.04%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '25.93'
$ws.Cells.Item(29, 5).Value = '  +2.58%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +2.11%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.37%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '10.28'
$ws.Cells.Item(32, 5).Value = '  +5.31%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '34.63'
$ws.Cells.Item(33, 5).Value = '  +6.62%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +1.57%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '2.06'
$ws.Cells.Item(35, 5).Value = '  +1.79%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '0.0436'
$ws.Cells.Item(36, 5).Value = '  +0.34%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.01%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +5.81%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.117'
$ws.Cells.Item(39, 5).Value = '  +2.47%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '1.84'
$ws.Cells.Item(40, 5).Value = '  +4.61%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '16.59'
$ws.Cells.Item(41, 5).Value = '  +3.78%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +4.37%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '125.25'
$ws.Cells.Item(43, 5).Value = '  +5.24%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '21.62'
$ws.Cells.Item(44, 5).Value = '  +3.87%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +11.70%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '2.02'
$ws.Cells.Item(46, 5).Value = '  +0.05%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '2.37'
$ws.Cells.Item(47, 5).Value = '  +4.07%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.23%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '2.032.17'
$ws.Cells.Item(49, 5).Value = '  +3.44%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '0.0323'
$ws.Cells.Item(50, 5).Value = '  +0.97%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '0.515'
$ws.Cells.Item(51, 5).Value = '  +13.45%  '

# Restore default cell style on column D (keeps text type) so the only
# differences from the source are the values themselves.
$ws.Range("D2:D51").Style = "Normal"
